# Refresh the cryptocurrency price/volume snapshot to the latest pull.
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.224.54"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.859.97"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.47"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6989"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07828"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3118"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.03"
$ws.Range("E10").Value = "  -3.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07799"
$ws.Range("E11").Value = "  -4.07%  "
$ws.Range("D12").Value = "1.863.54"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.134"
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.06"
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6935"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.565"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008535"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "29.259.77"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.37"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").Value = "2.114.80"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.93"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.64"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.912"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.60"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.580"
$ws.Range("E29").Value = "  +4.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.281"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.254"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05236"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.878"
$ws.Range("E35").Value = "  -3.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01855"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("D39").Value = "1.243.87"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.738"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9023"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.10"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.891"
$ws.Range("E43").Value = "  -7.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.70"
$ws.Range("E45").Value = "  -7.48%  "
$ws.Range("D46").Value = "2.012.53"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000124"
$ws.Range("E47").Value = "  -4.52%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5186"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.547"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.768"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4261"
$ws.Range("E51").Value = "  -1.90%  "
